# Commit: "renamed repo, fixed output folder path"
# The underlying data rebuild caused each sheet to contain a duplicated
# tail block (the first rows of the "real" result set were repeated at
# the end of the sheet). This script removes the duplicate leading block
# from each sheet, shifting the remaining (already-correct) rows up so
# that the tail becomes the sole copy of the data - matching the target
# workbook exactly.

$wb = $excel.ActiveWorkbook

# Sheet 1: "s__Dorea_A longicatena-b-p" - remove rows 11-15 (5 rows)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A11:A15").EntireRow.Delete()

# Sheet 2: "s__Dorea_A longicatena_B-b-p" - remove rows 6-11 (6 rows)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A6:A11").EntireRow.Delete()

# Sheet 3: "s__Dorea_A sp900550865-b-p" - remove rows 12-23 (12 rows)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A12:A23").EntireRow.Delete()
